# DistrictHeating 'smart heating' functionality
# Add two new conversionAssets rows (DH_Heat_Pump / DH_Peak_Boiler) and
# update storageAssets formatting / capacity values.

$wb = $excel.ActiveWorkbook

$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsStorage    = $wb.Worksheets.Item("storageAssets")

# --- conversionAssets: new row 9 (DH_Heat_Pump) ---
$wsConversion.Range("A9").Value = 8
$wsConversion.Range("B9").Value = "DH_Heat_Pump"
$wsConversion.Range("C9").Value = "CONVERSION"
$wsConversion.Range("D9").Value = "HEAT_PUMP_GROUND"
$wsConversion.Range("E9").Value = 100
$wsConversion.Range("F9").Value = 0
$wsConversion.Range("G9").Value = 0.65

# --- conversionAssets: new row 10 (DH_Peak_Boiler) ---
$wsConversion.Range("A10").Value = 9
$wsConversion.Range("B10").Value = "DH_Peak_Boiler"
$wsConversion.Range("C10").Value = "CONVERSION"
$wsConversion.Range("D10").Value = "BOILER"
$wsConversion.Range("E10").Value = 300
$wsConversion.Range("F10").Value = 297
$wsConversion.Range("G10").Value = 0.99

# --- storageAssets: give L3/L4 the same scientific-notation style already used by L5:L9 ---
$wsStorage.Range("L3").NumberFormat = "0.00E+00"
$wsStorage.Range("L4").NumberFormat = "0.00E+00"

# --- storageAssets: bump max capacity of last storage asset ---
$wsStorage.Range("L11").Value = 1000000000

# --- restore on-screen selections (storageAssets stays the active tab/sheet) ---
$wsConversion.Activate() | Out-Null
$wsConversion.Range("E10").Select() | Out-Null

$wsStorage.Activate() | Out-Null
$wsStorage.Range("L12").Select() | Out-Null
